$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 58 (shifting the old row 58 "**" marker to row 59)
$ws.Rows.Item(58).Insert()

$rng = $ws.Range("A58:G58")
$rng.NumberFormat = "@"

$ws.Cells.Item(58, 1).Value = "05/03/2018"
$ws.Cells.Item(58, 2).Value = "3070"
$ws.Cells.Item(58, 3).Value = 3
$ws.Cells.Item(58, 4).Value = "test"
$ws.Cells.Item(58, 5).Value = "1"
$ws.Cells.Item(58, 6).Value = "80000571"
$ws.Cells.Item(58, 7).Value = "245"
